$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": statuses for the two handed-back files have moved on.
# - 28c8dcfe-...: now "In Translation" (was "Handed back: in sync with en-US")
# - ca40efca-...: now "Ready for handoff" (was "Handed back: in sync with en-US")
# The "Latest HO Xliff Generate Date" for both rows also refreshed.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-12-16 09:38:32"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-12-16 09:38:32"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------------
# Sheet "zh-cn": per-language detail rows mirror the same status change, plus
# refreshed handoff datetime and a new handback-version warning message in
# the "Error Detail" column.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2016-12-16 09:38:18"
$wsZhCn.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/28c8dcfe-f72d-4a1e-8572-60038800e9d0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/655f5081a82cf01c4f1211e7165e0a2881104e7c/e2e/28c8dcfe-f72d-4a1e-8572-60038800e9d0.md."

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-12-16 09:38:18"
$wsZhCn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/ca40efca-d1b0-4c5f-af3c-95fc671a0aee.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/655f5081a82cf01c4f1211e7165e0a2881104e7c/e2e/ca40efca-d1b0-4c5f-af3c-95fc671a0aee.md."

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(18).ColumnWidth = 40

# ---------------------------------------------------------------------------
# Sheet "de-de": same as zh-cn, but keeps its own handoff datetime (unchanged)
# and xlf file names.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2016-12-16 09:38:32"
$wsDeDe.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/28c8dcfe-f72d-4a1e-8572-60038800e9d0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/655f5081a82cf01c4f1211e7165e0a2881104e7c/e2e/28c8dcfe-f72d-4a1e-8572-60038800e9d0.md."

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-12-16 09:38:32"
$wsDeDe.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/ca40efca-d1b0-4c5f-af3c-95fc671a0aee.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/655f5081a82cf01c4f1211e7165e0a2881104e7c/e2e/ca40efca-d1b0-4c5f-af3c-95fc671a0aee.md."

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(18).ColumnWidth = 40
